$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (pushes the existing row 12 -> 13 and row 13 -> 14,
# carrying their data/formatting down with them).
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 44474
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = 100112013
$ws.Range("G12").Value = "Alcachofa"
$ws.Range("H12").Value = "Española"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = "$/caja 30 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 333
$ws.Range("Q12").Value = 30
$ws.Range("R12").Value = "Hortaliza"

# Keep the date column formatted consistently with the rest of the sheet.
$ws.Range("D12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
